$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 updates
$ws.Range("M12").Value = 1.13
$ws.Range("N12").Value = 6
$ws.Range("AR12").Value = 2.1
$ws.Range("AS12").Value = 1.78

# Row 31 updates
$ws.Range("G31").Value = 2.4
$ws.Range("H31").Value = 3.2
$ws.Range("I31").Value = 3.1
$ws.Range("J31").Value = 3.1
$ws.Range("L31").Value = 3.75
$ws.Range("O31").Value = 1.36
$ws.Range("P31").Value = 3
$ws.Range("U31").Value = 1.44
$ws.Range("V31").Value = 2.63
$ws.Range("W31").Value = 1.8
$ws.Range("X31").Value = 1.95
$ws.Range("Z31").Value = 11
$ws.Range("AA31").Value = 9.5
$ws.Range("AB31").Value = 23
$ws.Range("AD31").Value = 29
$ws.Range("AE31").Value = 8.5
$ws.Range("AJ31").Value = 9
$ws.Range("AK31").Value = 15
$ws.Range("AM31").Value = 34
$ws.Range("AN31").Value = 26

# Row 60 updates
$ws.Range("G60").Value = 2.5
$ws.Range("H60").Value = 2.52
$ws.Range("J60").Value = 3.2
$ws.Range("K60").Value = 1.78
$ws.Range("L60").Value = 4.1
$ws.Range("M60").Value = 1.12
$ws.Range("N60").Value = 6
$ws.Range("O60").Value = 1.52
$ws.Range("P60").Value = 2.2
$ws.Range("Q60").Value = 2.5
$ws.Range("R60").Value = 1.4
$ws.Range("S60").Value = 4.2
$ws.Range("U60").Value = 1.57
$ws.Range("V60").Value = 2.12
$ws.Range("W60").Value = 1.98
$ws.Range("X60").Value = 1.65
$ws.Range("Y60").Value = 6
$ws.Range("Z60").Value = 11.25
$ws.Range("AB60").Value = 30
$ws.Range("AC60").Value = 26
$ws.Range("AD60").Value = 40
$ws.Range("AE60").Value = 5.5
$ws.Range("AF60").Value = 5.1
$ws.Range("AG60").Value = 15.5
$ws.Range("AH60").Value = 100
$ws.Range("AJ60").Value = 7.5
$ws.Range("AK60").Value = 17
$ws.Range("AM60").Value = 55
$ws.Range("AN60").Value = 40
